# Single and multi corrector: strip the English abbreviation + numbering
# out of column A, mirror the bare abbreviation (or leave blank) in column B,
# and fix a handful of OCR-style typos in column C/D (mmo1/L -> mmol/L,
# umol/L -> μmol/L, mo1/L -> mol/L, 4.20 -> 4.2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - 钾离子（K）
$ws.Range("A2").Value = "钾离子"
$ws.Range("B2").Value = "K"
$ws.Range("C2").Value = "'4.2"
$ws.Range("D2").Value = "mmol/L"

# Row 3 - 钠离子（Na）
$ws.Range("A3").Value = "钠离子"
$ws.Range("B3").ClearContents()
$ws.Range("D3").Value = "mmol/L"

# Row 4 - 氯离子（Cl）
$ws.Range("A4").Value = "氯离子"
$ws.Range("B4").Value = "Cl"
$ws.Range("D4").Value = "mmol/L"

# Row 5 - 总二氧化碳（TCO2）
$ws.Range("A5").Value = "总二氧化碳"
$ws.Range("B5").ClearContents()

# Row 6 - 尿素（UREA）
$ws.Range("A6").Value = "尿素"
$ws.Range("B6").ClearContents()

# Row 7 - 肌酐（CR）
$ws.Range("A7").Value = "肌酐"
$ws.Range("B7").ClearContents()

# Row 8 - 尿酸（UA）
$ws.Range("A8").Value = "尿酸"
$ws.Range("B8").ClearContents()
$ws.Range("D8").Value = "μmol/L"

# Row 9 - 丙氨酸氨基转移酶（ALT）
$ws.Range("A9").Value = "丙氨酸氨基转移酶"
$ws.Range("B9").ClearContents()

# Row 10 - 天门冬氨酸氨基转移酶（AST）
$ws.Range("A10").Value = "天门冬氨酸氨基转移酶"
$ws.Range("B10").ClearContents()

# Row 11 - 转氨酶比值（AST/ALT）
$ws.Range("A11").Value = "转氨酶比值"
$ws.Range("B11").ClearContents()

# Row 12 - 总蛋白（TP）
$ws.Range("A12").Value = "总蛋白"
$ws.Range("B12").ClearContents()

# Row 13 - 白蛋白（ALB）
$ws.Range("A13").Value = "白蛋白"
$ws.Range("B13").ClearContents()

# Row 14 - 球蛋白（G）
$ws.Range("A14").Value = "球蛋白"
$ws.Range("B14").ClearContents()

# Row 15 - 白/球比值（A/G）
$ws.Range("A15").Value = "白/球比值"
$ws.Range("B15").ClearContents()

# Row 16 - 总胆红素（TBIL）
$ws.Range("A16").Value = "总胆红素"
$ws.Range("B16").ClearContents()

# Row 17 - 直接胆红素（DBIL）
$ws.Range("A17").Value = "直接胆红素"
$ws.Range("B17").ClearContents()
$ws.Range("D17").Value = "μmol/L"

# Row 18 - 间接胆红素（IBIL）
$ws.Range("A18").Value = "间接胆红素"
$ws.Range("B18").ClearContents()
$ws.Range("D18").Value = "mol/L"

# Row 19 - 葡萄糖（GLU）
$ws.Range("A19").Value = "葡萄糖"
$ws.Range("B19").ClearContents()

# Row 20 - 甘油三酯（TG）
$ws.Range("A20").Value = "甘油三酯"
$ws.Range("B20").ClearContents()

# Row 21 - 甘油三酯（TG）
$ws.Range("A21").Value = "甘油三酯"
$ws.Range("B21").ClearContents()

# Row 22 - 总胆固醇（CHOL）
$ws.Range("A22").Value = "总胆固醇"
$ws.Range("B22").ClearContents()

# Row 23 - 总胆固醇（CHOL）
$ws.Range("A23").Value = "总胆固醇"
$ws.Range("B23").ClearContents()

# Row 24 - 高密度脂蛋白胆固醇（HDL-C）
$ws.Range("A24").Value = "高密度脂蛋白胆固醇"
$ws.Range("B24").ClearContents()
$ws.Range("D24").Value = "mmol/L"

# Row 25 - 高密度脂蛋白胆固醇（HDL-C）
$ws.Range("A25").Value = "高密度脂蛋白胆固醇"
$ws.Range("B25").ClearContents()

# Row 26 - 低密度脂蛋白胆固醇（LDL-C）
$ws.Range("A26").Value = "低密度脂蛋白胆固醇"
$ws.Range("B26").ClearContents()

# Row 27 - 低密度脂蛋白胆固醇（LDL-C）
$ws.Range("A27").Value = "低密度脂蛋白胆固醇"
$ws.Range("B27").ClearContents()

# Row 28 - 极低密度脂蛋白胆固醇（VLDL-C）
$ws.Range("A28").Value = "极低密度脂蛋白胆固醇"
$ws.Range("B28").ClearContents()

# Row 29 - 极低密度脂蛋白胆固醇（VLDL-C）
$ws.Range("A29").Value = "极低密度脂蛋白胆固醇"
$ws.Range("B29").ClearContents()

# Row 30 - 非高密度脂蛋白胆固醇（nHDLC）
$ws.Range("A30").Value = "非高密度脂蛋白胆固醇"
$ws.Range("B30").ClearContents()

# Row 31 - 非高密度脂蛋白胆固醇（nHDLC）
$ws.Range("A31").Value = "非高密度脂蛋白胆固醇"
$ws.Range("B31").ClearContents()
